$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Insert new row 29 (numArr code) after existing row 28, and
#    new row 31 (calculator function) after row 30 (old row 29),
#    but set the *content* in the same order the new shared
#    strings appear in the original commit (calc, ascii x2, numArr)
#    so the sharedStrings table matches byte-for-byte.
# ---------------------------------------------------------------
$ws.Rows.Item(29).Insert()
$ws.Rows.Item(31).Insert()

# calculator function (ends up in row 31)
$calcCode = "function calculaton(num1, num2, opr){`n    let result = {`n        'multiply': (a, b) => a * b,`n        'divide': (a, b) => a / b,`n        'add': (a, b) => a + b,`n        'subtract': (a, b) => a - b`n    }`n    console.log(result[opr](num1, num2));`n}"
$ws.Range("A31").Value2 = $calcCode
$ws.Range("B31").Value2 = "Калкоратор чрез речник без if - else проверки с използването на arow function"
$ws.Rows.Item(31).RowHeight = 148.8
$ws.Range("B31").HorizontalAlignment = -4108
$ws.Range("B31").VerticalAlignment = -4108
$ws.Range("B31").WrapText = $true
$ws.Range("B31").Font.Size = 12

# ASCII helpers (fill previously-empty rows 38 and 39 - no insert
# needed, these row numbers were unused blank space in the sheet,
# and simply take the column default formatting.)
$ws.Range("A38").Value2 = "chr1.charCodeAt(0);"
$ws.Range("B38").Value2 = "Връща ASCII стойността на символа."

$ws.Range("A39").Value2 = "String.fromCharCode(i);"
$ws.Range("B39").Value2 = "Връща символ от ASCII таблицата по зададено число."

# number-to-digit-array helper (ends up in row 29)
$ws.Range("A29").Value2 = "let numArr = String(num).split(`"`").map((num) => Number(num))"
$ws.Range("B29").Value2 = "Превръща число в масив от числа . 12345 - [1, 2, 3, 4, 5]"
$ws.Rows.Item(29).RowHeight = 18

# ---------------------------------------------------------------
# 4) Column / view adjustments
# ---------------------------------------------------------------
# Target stored width is 70.109375 characters; the engine quantizes
# ColumnWidth to 1/6-character steps, so 69.35 is the closest input
# that lands on the nearest achievable stored width (70.1666...).
$ws.Columns.Item(1).ColumnWidth = 69.35

$ws.Application.ActiveWindow.ScrollRow = 20
$ws.Range("B31").Select()
